$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Remove all existing hyperlinks (and the relationships backing
#    them) so we can rebuild the ref -> target mapping cleanly.
# ------------------------------------------------------------------
$ws.Hyperlinks.Delete()

# ------------------------------------------------------------------
# 2. Clear the two cells that disappear entirely in the new layout
#    (thumbnail moves out of H2, affiliation moves out of M2).
# ------------------------------------------------------------------
$ws.Range("H2").Clear()
$ws.Range("M2").Clear()

# ------------------------------------------------------------------
# 3. Header row (row 1) text updates.
# ------------------------------------------------------------------
$ws.Range("F1").Value = "機械可読ドキュメント"
$ws.Range("H1").Value = "年"
$ws.Range("I1").Value = "IIIFマニフェストURI"
$ws.Range("J1").Value = "viewingDirection"
$ws.Range("K1").Value = "帰属"
$ws.Range("M1").Value = "ソート用項目"
$ws.Range("N1").Value = "コレクション"
$ws.Range("O1").Value = "サムネイル"
$ws.Range("P1").Value = "ウェブサイトURL"

# ------------------------------------------------------------------
# 4. Plain (non-hyperlinked) value cells in row 2.
# ------------------------------------------------------------------
$ws.Range("K2").Value = "東京大学総合図書館 General Library in the University of Tokyo, JAPAN"
$ws.Range("N2").Value = "直江状"

# ------------------------------------------------------------------
# 5. Re-create hyperlinks in document order so the relationship ids
#    line up (rId1..rId10) the same way the source file uses them.
#    TextToDisplay pins each cell's visible text to the desired
#    final value in the same call.
# ------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("C2"), "https://www.lib.u-tokyo.ac.jp/ja/library/general/reuse", [Type]::Missing, [Type]::Missing, "https://www.lib.u-tokyo.ac.jp/ja/library/general/reuse")
$ws.Hyperlinks.Add($ws.Range("D2"), "https://iiif.dl.itc.u-tokyo.ac.jp/repo/s/naoe/document/daaf81b0-f1a7-435c-9216-a8fb56390887", [Type]::Missing, [Type]::Missing, "https://iiif.dl.itc.u-tokyo.ac.jp/repo/s/naoe/document/daaf81b0-f1a7-435c-9216-a8fb56390887")
$ws.Hyperlinks.Add($ws.Range("F2"), "https://iiif.dl.itc.u-tokyo.ac.jp/repo/api/items/13", [Type]::Missing, [Type]::Missing, "https://iiif.dl.itc.u-tokyo.ac.jp/repo/api/items/13")
$ws.Hyperlinks.Add($ws.Range("I2"), "https://iiif.dl.itc.u-tokyo.ac.jp/repo/iiif-img/85/full/200,151/0/default.jpg", [Type]::Missing, [Type]::Missing, "https://iiif.dl.itc.u-tokyo.ac.jp/repo/iiif/13/manifest")
$ws.Hyperlinks.Add($ws.Range("J2"), "https://iiif.dl.itc.u-tokyo.ac.jp/repo/iiif/13/manifest", "rightToLeftDirection", [Type]::Missing, "http://iiif.io/api/presentation/2#rightToLeftDirection")
$ws.Hyperlinks.Add($ws.Range("O2"), "http://iiif.io/api/presentation/2", [Type]::Missing, [Type]::Missing, "https://iiif.dl.itc.u-tokyo.ac.jp/repo/iiif-img/85/full/200,151/0/default.jpg")
$ws.Hyperlinks.Add($ws.Range("P2"), "https://iiif.dl.itc.u-tokyo.ac.jp/repo/s/naoe/", [Type]::Missing, [Type]::Missing, "https://iiif.dl.itc.u-tokyo.ac.jp/repo/s/naoe/")
$ws.Hyperlinks.Add($ws.Range("Q2"), "https://iiif.dl.itc.u-tokyo.ac.jp/omekac/oa/collections/8/manifest.json", [Type]::Missing, [Type]::Missing, "https://iiif.dl.itc.u-tokyo.ac.jp/omekac/oa/collections/8/manifest.json")
$ws.Hyperlinks.Add($ws.Range("R2"), "http://tapasproject.org/tapas-commons/files/直江状", [Type]::Missing, [Type]::Missing, "http://tapasproject.org/tapas-commons/files/直江状")
$ws.Hyperlinks.Add($ws.Range("S2"), "https://iiif.dl.itc.u-tokyo.ac.jp/api/iiif-search/kPzFpI4mtex7HdRmrZL1ew9r7OCgdDPvNX2g0njpVtAV%EF%BC%8BWUgecS%EF%BC%8BSsVMvlKTXaNVmk9OVUlQkEsXFV%EF%BC%8B86MGLyr5YB2lgAB9MM6QnC63BVFQ%3D", [Type]::Missing, [Type]::Missing, "https://iiif.dl.itc.u-tokyo.ac.jp/api/iiif-search/kPzFpI4mtex7HdRmrZL1ew9r7OCgdDPvNX2g0njpVtAV%EF%BC%8BWUgecS%EF%BC%8BSsVMvlKTXaNVmk9OVUlQkEsXFV%EF%BC%8B86MGLyr5YB2lgAB9MM6QnC63BVFQ%3D")

# ------------------------------------------------------------------
# 6. Hyperlinks.Add stamps a brand-new "Hyperlink" style (theme
#    colour) on each touched cell; restore the workbook's original
#    custom hyperlink font (underline + RGB 0000FF) so the cells
#    collapse back onto the pre-existing style used throughout the
#    sheet.
# ------------------------------------------------------------------
foreach ($addr in @("C2","D2","F2","I2","J2","O2","P2","Q2","R2","S2")) {
    $rng = $ws.Range($addr)
    $rng.Font.Underline = $true
    $rng.Font.Color = 16711680
}
